$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 4.91503868274033
$ws.Cells.Item(2, 4).Value = 8.431355421129439
$ws.Cells.Item(2, 5).Value = 16.32120603405771
$ws.Cells.Item(2, 6).Value = 47.03836406923052
$ws.Cells.Item(2, 7).Value = 3.670819438127647
$ws.Cells.Item(2, 11).Value = 18.38214355479963
$ws.Cells.Item(3, 3).Value = 4.751603684059941
$ws.Cells.Item(3, 4).Value = 8.303994852852774
$ws.Cells.Item(3, 5).Value = 15.39775253343196
$ws.Cells.Item(3, 6).Value = 45.75581943335377
$ws.Cells.Item(3, 7).Value = 3.676673723525904
$ws.Cells.Item(3, 11).Value = 17.9275243313831
$ws.Cells.Item(4, 3).Value = 4.650282341326423
$ws.Cells.Item(4, 4).Value = 8.225840486797102
$ws.Cells.Item(4, 5).Value = 14.80847425741368
$ws.Cells.Item(4, 6).Value = 44.96147644748883
$ws.Cells.Item(4, 7).Value = 3.680437462719171
$ws.Cells.Item(4, 11).Value = 17.65165409275166
$ws.Cells.Item(5, 3).Value = 4.608820457233893
$ws.Cells.Item(5, 4).Value = 8.194022593157518
$ws.Cells.Item(5, 5).Value = 14.56301185418733
$ws.Cells.Item(5, 6).Value = 44.63647258934079
$ws.Cells.Item(5, 7).Value = 3.682014011228191
$ws.Cells.Item(5, 11).Value = 17.54025764445091
$ws.Cells.Item(6, 3).Value = 4.601927506422632
$ws.Cells.Item(6, 4).Value = 8.188741661046778
$ws.Cells.Item(6, 5).Value = 14.52194035389332
$ws.Cells.Item(6, 6).Value = 44.58243994344635
$ws.Cells.Item(6, 7).Value = 3.682278387913593
$ws.Cells.Item(6, 11).Value = 17.52182806753881
$ws.Cells.Item(7, 3).Value = 4.649723772482862
$ws.Cells.Item(7, 4).Value = 8.225411230375277
$ws.Cells.Item(7, 5).Value = 14.80518503190238
$ws.Cells.Item(7, 6).Value = 44.95709803579393
$ws.Cells.Item(7, 7).Value = 3.680458550998166
$ws.Cells.Item(7, 11).Value = 17.65014735140937
$ws.Cells.Item(8, 3).Value = 4.858929443519326
$ws.Cells.Item(8, 4).Value = 8.38744364215739
$ws.Cells.Item(8, 5).Value = 16.00757100722909
$ws.Cells.Item(8, 6).Value = 46.59781212144669
$ws.Cells.Item(8, 7).Value = 3.672803033779343
$ws.Cells.Item(8, 11).Value = 18.22483218937939
$ws.Cells.Item(9, 3).Value = 5.258580251619456
$ws.Cells.Item(9, 4).Value = 8.704700534695242
$ws.Cells.Item(9, 5).Value = 18.21314253706674
$ws.Cells.Item(9, 6).Value = 49.74341967607151
$ws.Cells.Item(9, 7).Value = 3.659121444324172
$ws.Cells.Item(9, 11).Value = 19.36911012000945
$ws.Cells.Item(10, 3).Value = 5.542251486543845
$ws.Cells.Item(10, 4).Value = 8.936358318052488
$ws.Cells.Item(10, 5).Value = 19.8650251140621
$ws.Cells.Item(10, 6).Value = 51.98906720236063
$ws.Cells.Item(10, 7).Value = 3.649864626687362
$ws.Cells.Item(10, 11).Value = 20.2093545199134
$ws.Cells.Item(11, 3).Value = 5.668512670147665
$ws.Cells.Item(11, 4).Value = 9.04119765235502
$ws.Cells.Item(11, 5).Value = 20.57613065875274
$ws.Cells.Item(11, 6).Value = 52.99240436809512
$ws.Cells.Item(11, 7).Value = 3.645822553192565
$ws.Cells.Item(11, 11).Value = 20.58936138233784
$ws.Cells.Item(12, 3).Value = 5.715879428359457
$ws.Cells.Item(12, 4).Value = 9.080799962636652
$ws.Cells.Item(12, 5).Value = 20.83965277867336
$ws.Cells.Item(12, 6).Value = 53.36941990072702
$ws.Cells.Item(12, 7).Value = 3.644315930736049
$ws.Cells.Item(12, 11).Value = 20.7327770758288
$ws.Cells.Item(13, 3).Value = 5.705698638392763
$ws.Cells.Item(13, 4).Value = 9.072275570628467
$ws.Cells.Item(13, 5).Value = 20.78315397868444
$ws.Cells.Item(13, 6).Value = 53.28835760578055
$ws.Cells.Item(13, 7).Value = 3.644639344644899
$ws.Cells.Item(13, 11).Value = 20.70191398808173
$ws.Cells.Item(14, 3).Value = 5.672418734992571
$ws.Cells.Item(14, 4).Value = 9.044457836228318
$ws.Cells.Item(14, 5).Value = 20.5979259348419
$ws.Cells.Item(14, 6).Value = 53.02348197231183
$ws.Cells.Item(14, 7).Value = 3.645698122523621
$ws.Cells.Item(14, 11).Value = 20.60117089110901
$ws.Cells.Item(15, 3).Value = 5.651974538542479
$ws.Cells.Item(15, 4).Value = 9.027405315242685
$ws.Cells.Item(15, 5).Value = 20.48371982026693
$ws.Cells.Item(15, 6).Value = 52.86084825480727
$ws.Cells.Item(15, 7).Value = 3.646349775158401
$ws.Cells.Item(15, 11).Value = 20.53939503792302
$ws.Cells.Item(16, 3).Value = 5.533940024560027
$ws.Cells.Item(16, 4).Value = 8.929494173847173
$ws.Cells.Item(16, 5).Value = 19.81774456872767
$ws.Cells.Item(16, 6).Value = 51.92310330769706
$ws.Cells.Item(16, 7).Value = 3.650132161836786
$ws.Cells.Item(16, 11).Value = 20.18446105026619
$ws.Cells.Item(17, 3).Value = 5.460783016452749
$ws.Cells.Item(17, 4).Value = 8.869275403159387
$ws.Cells.Item(17, 5).Value = 19.39888802430311
$ws.Cells.Item(17, 6).Value = 51.34293392770114
$ws.Cells.Item(17, 7).Value = 3.652495605400301
$ws.Cells.Item(17, 11).Value = 19.9660296934871
$ws.Cells.Item(18, 3).Value = 5.418445541145674
$ws.Cells.Item(18, 4).Value = 8.8345892198179
$ws.Cells.Item(18, 5).Value = 19.1541728429543
$ws.Cells.Item(18, 6).Value = 51.00753755292568
$ws.Cells.Item(18, 7).Value = 3.65387090921125
$ws.Cells.Item(18, 11).Value = 19.84019346950317
$ws.Cells.Item(19, 3).Value = 5.40406770275662
$ws.Cells.Item(19, 4).Value = 8.822837115119365
$ws.Cells.Item(19, 5).Value = 19.07066252029075
$ws.Cells.Item(19, 6).Value = 50.89369621216496
$ws.Cells.Item(19, 7).Value = 3.654339304629665
$ws.Cells.Item(19, 11).Value = 19.79755830740032
$ws.Cells.Item(20, 3).Value = 5.468597908251446
$ws.Cells.Item(20, 4).Value = 8.875691109180126
$ws.Cells.Item(20, 5).Value = 19.44386898781196
$ws.Cells.Item(20, 6).Value = 51.40487202354698
$ws.Cells.Item(20, 7).Value = 3.652242367710552
$ws.Cells.Item(20, 11).Value = 19.98930404554962
$ws.Cells.Item(21, 3).Value = 5.682206280372289
$ws.Cells.Item(21, 4).Value = 9.052631397291929
$ws.Cells.Item(21, 5).Value = 20.65248782938567
$ws.Cells.Item(21, 6).Value = 53.10136409750766
$ws.Cells.Item(21, 7).Value = 3.645386483937823
$ws.Cells.Item(21, 11).Value = 20.63077598637985
$ws.Cells.Item(22, 3).Value = 5.819196796718298
$ws.Cells.Item(22, 4).Value = 9.167691625065979
$ws.Cells.Item(22, 5).Value = 21.40884756997446
$ws.Cells.Item(22, 6).Value = 54.19293907181471
$ws.Cells.Item(22, 7).Value = 3.641045667506076
$ws.Cells.Item(22, 11).Value = 21.04712235068359
$ws.Cells.Item(23, 3).Value = 5.746335364158111
$ws.Cells.Item(23, 4).Value = 9.106341383913218
$ws.Cells.Item(23, 5).Value = 21.00821895113341
$ws.Cells.Item(23, 6).Value = 53.61201120757797
$ws.Cells.Item(23, 7).Value = 3.643349730190228
$ws.Cells.Item(23, 11).Value = 20.825226377287
$ws.Cells.Item(24, 3).Value = 5.465065659739312
$ws.Cells.Item(24, 4).Value = 8.872790771325123
$ws.Cells.Item(24, 5).Value = 19.42354527306953
$ws.Cells.Item(24, 6).Value = 51.37687554691497
$ws.Cells.Item(24, 7).Value = 3.652356804967226
$ws.Cells.Item(24, 11).Value = 19.97878249923743
$ws.Cells.Item(25, 3).Value = 5.151980689174306
$ws.Cells.Item(25, 4).Value = 8.619043711159742
$ws.Cells.Item(25, 5).Value = 17.61277334444323
$ws.Cells.Item(25, 6).Value = 48.90247185977486
$ws.Cells.Item(25, 7).Value = 3.662681874125008
$ws.Cells.Item(25, 11).Value = 19.05886960360095
